$wb = $excel.ActiveWorkbook

# Add the new "2025" worksheet after the last existing sheet ("2024")
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "2025"

# Populate the data for the new sheet (same layout/labels as the other year tabs)
$newSheet.Range("A1").Value = "Grupos de Edad"
$newSheet.Range("B1").Value = "Mujeres"
$newSheet.Range("C1").Value = "Hombres"

$newSheet.Range("A2").Value = "Menores"
$newSheet.Range("B2").Value = 16
$newSheet.Range("C2").Value = 165

$newSheet.Range("A3").Value = "Mayores"
$newSheet.Range("B3").Value = 165
$newSheet.Range("C3").Value = 630

$newSheet.Range("A4").Value = "Desconocida"
$newSheet.Range("B4").Value = 67
$newSheet.Range("C4").Value = 371

# Match column A width to content / the saved width from the authored workbook
$newSheet.Columns.Item(1).ColumnWidth = 16.5

# Select the whole data range on the new sheet (matches saved selection)
$null = $newSheet.Range("A1:C4").Select()

# The previously-active tab ("2024") remains the active/selected tab
$wb.Worksheets.Item("2024").Activate()
